# Atualiza planilhas mensais (run 2025-09-18)
# Insere 4 novas linhas de "Linguiça Toscana Swift 700 g" logo acima da
# linha "Lingüiça Toscana Grossa Auora Aprox. 700g" (antiga linha 71),
# empurrando as demais linhas para baixo (71-106 -> 75-110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere 4 linhas em branco acima da linha 71 (shift down).
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(71).Insert()
}

# Preenche as 4 novas linhas com os dados do produto "Linguiça Toscana Swift 700 g".
for ($r = 71; $r -le 74; $r++) {
    $ws.Cells.Item($r, 1).Value = "Linguiça Toscana Swift 700 g"
    $ws.Cells.Item($r, 2).Value = 14.9
    $ws.Cells.Item($r, 3).Value = 14.9
    $ws.Cells.Item($r, 4).Value = 14.9
}
